$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits at the very start of the
#    document (around the title "Java OOP Exam - 10 April 2021").
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) In the "FreshwaterAquarium" section ("Has 50 capacity" / "The
#    constructor should take the following values upon initialization:"),
#    bold the word "should" and wrap "should " in a brand new "_GoBack"
#    bookmark, splitting the run into three runs:
#      "should" (bold) + " " + "take the following values upon initialization:"
# ---------------------------------------------------------------------------

# Scope the search using the unique preceding text "Has 50 capacity" so we
# land on the one, specific occurrence of the target phrase (it repeats
# several times throughout the document).
$scope = $d.Range(0, $d.Content.End)
$scope.Find.Execute("Has 50 capacity", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($scope.End, $d.Content.End)
$target.Find.Execute("should take the following values upon initialization:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$phraseStart = $target.Start

# Shrink the run down to just "should" and make it bold.
$target.Text = "should"
$target.Bold = 1

# Re-insert the separating space as its own (non-bold) run.
$target.InsertAfter(" ")

# Append the remainder of the sentence as a further separate run.
$rest = $d.Range($target.End, $target.End)
$rest.InsertAfter("take the following values upon initialization:")

# Wrap "should " (the word plus the following space) in a new "_GoBack"
# bookmark.
$bmRange = $d.Range($phraseStart, $phraseStart + 7)
$d.Bookmarks.Add("_GoBack", $bmRange)
